$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.440.04'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.36%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.089.01'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +1.55%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '521.80'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.51'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.36%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.090.19'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +1.64%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.438'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.20'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.68%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.107'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.57%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.388'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +3.54%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.622.77'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +1.63%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.132'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.50%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.50'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.89%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000163'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.73%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '57.570.38'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.15%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.092.12'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +1.93%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.07'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.42%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.72'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.93'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '339.69'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +3.10%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.511'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.93%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '66.94'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +3.19%  '
$ws.Range('E26').Value = '  -0.18%  '
$ws.Range('D28').Value = '0.0₃0911'
$ws.Range('E28').Value = '  +1.81%  '
$ws.Range('B29').Value = 'USDe'
$ws.Range('C29').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.44'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.19'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.86'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +3.36%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.92'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +1.60%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.18'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '155.35'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.59%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.59'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.12'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +2.25%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '26.74'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.21%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.24'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0660'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.99%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.96'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +2.11%  '
$ws.Range('B42').Value = 'Mantle'
$ws.Range('C42').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.681'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +4.54%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.51'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +9.47%  '
$ws.Range('B44').Value = 'RenzoRestakedETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.129.83'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +1.61%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '36.68'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.38%  '
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.282.93'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0260'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.979'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +5.04%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.45'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.91%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.02'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +2.04%  '
